$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions (I1, J1) - copy formatting (style) from existing header cell H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data column I (I2:I8)
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 5
$ws.Range("I8").Value = 4

# Data column J (J2:J8)
$ws.Range("J2").Value = 1
$ws.Range("J3").Value = 6
$ws.Range("J4").Value = 7
$ws.Range("J5").Value = 5
$ws.Range("J6").Value = 5
$ws.Range("J7").Value = 6
$ws.Range("J8").Value = 5
